# Auto-generated edit script: apply latest cryptos snapshot (Price / Volume(1h) columns)
# Source: GitHub Actions "Updated cryptos list" commit — updates columns D (Price) and E (Volume 1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '42.589.83'
$ws.Range("E2").Value = '  -1.05%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '2.284.10'
# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.06%  '
# Row 5: BNB
$c = $ws.Range("D5")
$c.NumberFormat = "@"   # force text so Excel does not coerce "303.87" into a Number
$c.Value = '303.87'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E5").Value = '  +0.98%  '
# Row 6: Solana
$c = $ws.Range("D6")
$c.NumberFormat = "@"   # force text so Excel does not coerce "95.71" into a Number
$c.Value = '95.71'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E6").Value = '  -2.63%  '
# Row 8: USDC
$ws.Range("E8").Value = '  +0.02%  '
# Row 10: Avalanche
$c = $ws.Range("D10")
$c.NumberFormat = "@"   # force text so Excel does not coerce "34.68" into a Number
$c.Value = '34.68'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E10").Value = '  -4.15%  '
# Row 11: Dogecoin
$ws.Range("E11").Value = '  -1.21%  '
# Row 12: TRON
$ws.Range("E12").Value = '  +1.73%  '
# Row 13: Chainlink
$c = $ws.Range("D13")
$c.NumberFormat = "@"   # force text so Excel does not coerce "18.03" into a Number
$c.Value = '18.03'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E13").Value = '  +0.97%  '
# Row 14: Polkadot
$c = $ws.Range("D14")
$c.NumberFormat = "@"   # force text so Excel does not coerce "6.79" into a Number
$c.Value = '6.79'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E14").Value = '  -0.49%  '
# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '2.640.19'
$ws.Range("E15").Value = '  -0.86%  '
# Row 16: WrappedEther
$ws.Range("D16").Value = '2.267.94'
$ws.Range("E16").Value = '  -1.64%  '
# Row 17: Polygon
$ws.Range("E17").Value = '  -1.25%  '
# Row 18: WrappedBTC
$ws.Range("D18").Value = '42.508.64'
$ws.Range("E18").Value = '  -1.15%  '
# Row 19: InternetComputer(DFINITY)
$c = $ws.Range("D19")
$c.NumberFormat = "@"   # force text so Excel does not coerce "12.98" into a Number
$c.Value = '12.98'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E19").Value = '  +2.25%  '
# Row 20: ShibaInu
$ws.Range("E20").Value = '  -1.91%  '
# Row 21: Uniswap
$ws.Range("E21").Value = '  -2.66%  '
# Row 22: Litecoin
$c = $ws.Range("D22")
$c.NumberFormat = "@"   # force text so Excel does not coerce "67.13" into a Number
$c.Value = '67.13'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E22").Value = '  -1.73%  '
# Row 23: BitcoinCash
$c = $ws.Range("D23")
$c.NumberFormat = "@"   # force text so Excel does not coerce "235.69" into a Number
$c.Value = '235.69'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E23").Value = '  -2.82%  '
# Row 24: ImmutableX
$c = $ws.Range("D24")
$c.NumberFormat = "@"   # force text so Excel does not coerce "2.13" into a Number
$c.Value = '2.13'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E24").Value = '  -1.25%  '
# Row 25: Dai
$ws.Range("E25").Value = '  +0.19%  '
# Row 26: PancakeSwap
$ws.Range("E26").Value = '  -1.26%  '
# Row 27: EthereumClassic
$c = $ws.Range("D27")
$c.NumberFormat = "@"   # force text so Excel does not coerce "24.62" into a Number
$c.Value = '24.62'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E27").Value = '  -2.65%  '
# Row 28: Toncoin
$ws.Range("E28").Value = '  +16.98%  '
# Row 29: Monero
$c = $ws.Range("D29")
$c.NumberFormat = "@"   # force text so Excel does not coerce "166.38" into a Number
$c.Value = '166.38'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E29").Value = '  -0.38%  '
# Row 30: Cosmos
$c = $ws.Range("D30")
$c.NumberFormat = "@"   # force text so Excel does not coerce "8.96" into a Number
$c.Value = '8.96'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E30").Value = '  -1.49%  '
# Row 31: InjectiveProtocol
$c = $ws.Range("D31")
$c.NumberFormat = "@"   # force text so Excel does not coerce "32.64" into a Number
$c.Value = '32.64'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
# Row 32: FirstDigitalUSD
$ws.Range("E32").Value = '  +0.04%  '
# Row 33: Celestia
$c = $ws.Range("D33")
$c.NumberFormat = "@"   # force text so Excel does not coerce "17.85" into a Number
$c.Value = '17.85'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E33").Value = '  +0.65%  '
# Row 34: Filecoin
$ws.Range("E34").Value = '  -1.73%  '
# Row 35: RenderToken
$ws.Range("E35").Value = '  -6.40%  '
# Row 36: WEMIXToken
$c = $ws.Range("D36")
$c.NumberFormat = "@"   # force text so Excel does not coerce "2.36" into a Number
$c.Value = '2.36'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E36").Value = '  -2.31%  '
# Row 37: Hedera
$ws.Range("E37").Value = '  -1.00%  '
# Row 39: ARBITRUM
$ws.Range("E39").Value = '  -2.33%  '
# Row 40: Stellar
$ws.Range("E40").Value = '  -1.98%  '
# Row 41: LidoDAOToken
$c = $ws.Range("D41")
$c.NumberFormat = "@"   # force text so Excel does not coerce "2.67" into a Number
$c.Value = '2.67'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E41").Value = '  -3.86%  '
# Row 42: Maker
$ws.Range("D42").Value = '1.991.71'
$ws.Range("E42").Value = '  -0.73%  '
# Row 43: VeChain
$ws.Range("E43").Value = '  -3.69%  '
# Row 44: FraxShare
$c = $ws.Range("D44")
$c.NumberFormat = "@"   # force text so Excel does not coerce "10.19" into a Number
$c.Value = '10.19'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E44").Value = '  +0.02%  '
# Row 45: EnergySwap
$c = $ws.Range("D45")
$c.NumberFormat = "@"   # force text so Excel does not coerce "18.17" into a Number
$c.Value = '18.17'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E45").Value = '  +3.32%  '
# Row 46: ApeXProtocol
$c = $ws.Range("D46")
$c.NumberFormat = "@"   # force text so Excel does not coerce "2.03" into a Number
$c.Value = '2.03'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E46").Value = '  -6.29%  '
# Row 47: NEARProtocol
$c = $ws.Range("D47")
$c.NumberFormat = "@"   # force text so Excel does not coerce "2.75" into a Number
$c.Value = '2.75'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E47").Value = '  -2.20%  '
# Row 48: HuobiToken
$c = $ws.Range("D48")
$c.NumberFormat = "@"   # force text so Excel does not coerce "2.92" into a Number
$c.Value = '2.92'
$c.Style = "Normal"     # drop the temporary Text format, keep default styling
$ws.Range("E48").Value = '  -4.89%  '
# Row 49: MultiversX
$ws.Range("E49").Value = '  -0.25%  '
# Row 50: RocketPoolETH
$ws.Range("D50").Value = '2.505.28'
$ws.Range("E50").Value = '  -0.94%  '
# Row 51: TrustWalletToken
$ws.Range("E51").Value = '  +0.43%  '
